# "Generate Report for Handoff"
#
# The localization-status report previously reflected a finished handback
# cycle ("Handed back: in sync with en-US" with the handback timestamps).
# This run instead generates the report right after a fresh handoff, so:
#   - the Status cells flip to "Ready for handoff"
#   - the associated timestamps move to the handoff generation time
#   - the now-shorter status text no longer needs the wide status columns,
#     so those columns are narrowed to fit.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 20:57:58"

# Status columns (E & F) were sized for the old, longer status text;
# narrow them now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 20:57:52"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 20:57:58"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

Write-Host "Overview!E2 =" $wsOverview.Range("E2").Value()
Write-Host "Overview!F2 =" $wsOverview.Range("F2").Value()
Write-Host "Overview!G2 =" $wsOverview.Range("G2").Value()
Write-Host "zh-cn!C2 =" $wsZhCn.Range("C2").Value()
Write-Host "zh-cn!H2 =" $wsZhCn.Range("H2").Value()
Write-Host "de-de!C2 =" $wsDeDe.Range("C2").Value()
Write-Host "de-de!H2 =" $wsDeDe.Range("H2").Value()
